$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Text (shared-string) writes, in the exact order needed to reproduce
#     the target shared-strings table layout ---

# Row 12: Udarbejdelse af iterationsplan 2 / Project Manager
$ws.Range("A12").Value = "Udarbejdelse af iterationsplan 2"
$ws.Range("B12").Value = "Project Manager"

# Row 13: Kundemøde / Business-Process Analyst
$ws.Range("A13").Value = "Kundemøde"
$ws.Range("B13").Value = "Business-Process Analyst"

# Row 14: Virksomhedsbesøg / Any role
$ws.Range("A14").Value = "Virksomhedsbesøg"
$ws.Range("B14").Value = "Any role"

# Row 16 filled before row 15 (matches original authoring order)
$ws.Range("A16").Value = "Udarbejdelse af AD06"
$ws.Range("B16").Value = "Business-Process Analyst"

$ws.Range("A15").Value = "Udarbejdelse af AD05"
$ws.Range("B15").Value = "Business-Process Analyst"

# Row 17: Udarbejdelse af DD05 / Business-Process Analyst
$ws.Range("A17").Value = "Udarbejdelse af DD05"
$ws.Range("B17").Value = "Business-Process Analyst"

# Row 18: Udarbejdelse af DD06 / Business-Process Analyst
$ws.Range("A18").Value = "Udarbejdelse af DD06"
$ws.Range("B18").Value = "Business-Process Analyst"

# Row 19: Udarbejdelse af Mockup06b / Business-Process Analyst
$ws.Range("A19").Value = "Udarbejdelse af Mockup06b"
$ws.Range("B19").Value = "Business-Process Analyst"

# Row 20: Klasseundervisning i SD + DCD / Any role
$ws.Range("A20").Value = "Klasseundervisning i SD + DCD"
$ws.Range("B20").Value = "Any role"

# --- Numeric writes: date + start/end time fractions ---

$ws.Range("C12").Value = 43887
$ws.Range("D12").Value = 0.35416666666666669
$ws.Range("E12").Value = 0.41666666666666669

$ws.Range("C13").Value = 43887
$ws.Range("D13").Value = 0.43055555555555558
$ws.Range("E13").Value = 0.58333333333333337

$ws.Range("C14").Value = 43887
$ws.Range("D14").Value = 0.60416666666666663
$ws.Range("E14").Value = 0.66666666666666663

$ws.Range("C15").Value = 43888
$ws.Range("D15").Value = 0.35416666666666669
$ws.Range("E15").Value = 0.375

$ws.Range("C16").Value = 43888
$ws.Range("D16").Value = 0.375
$ws.Range("E16").Value = 0.39583333333333331

$ws.Range("C17").Value = 43888
$ws.Range("D17").Value = 0.39583333333333331
$ws.Range("E17").Value = 0.41666666666666669

$ws.Range("C18").Value = 43888
$ws.Range("D18").Value = 0.41666666666666669
$ws.Range("E18").Value = 0.4375

$ws.Range("C19").Value = 43888
$ws.Range("D19").Value = 0.4375
$ws.Range("E19").Value = 0.45833333333333331

$ws.Range("C20").Value = 43888
$ws.Range("D20").Value = 0.47916666666666669
$ws.Range("E20").Value = 0.66666666666666663

# Update the active selection to match the post-edit state
$ws.Activate() | Out-Null
$ws.Range("E20").Select() | Out-Null
